# Add team win/loss/tie record columns (AD, AE, AF) to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header styling (bold, bordered, centered) used by the rest of row 1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row shares the same 2016 team record: 78 wins, 83 losses, 1 tie
for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 30).Value = 78
    $ws.Cells.Item($row, 31).Value = 83
    $ws.Cells.Item($row, 32).Value = 1
}
